$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(239, 44313, 2, 32, 186.752261453166),
    @(240, 44314, 1, 27, 157.5722206011089),
    @(241, 44315, 6, 28, 163.4082287715203),
    @(242, 44316, 8, 33, 192.5882696235775),
    @(243, 44317, 4, 28, 163.4082287715203),
    @(244, 44318, 7, 33, 192.5882696235775)
)

# Copy the date-cell style from the last existing row (A238) so the new
# date cells (A239:A244) pick up the same number format / alignment / border
# instead of Excel fabricating a brand-new style entry.
$ws.Range("A238").Copy()

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("A$r").PasteSpecial(-4122)  # xlPasteFormats
}
